$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = -2.796573638916016
$ws.Cells.Item(2,4).Value = -10.2697286605835
$ws.Cells.Item(2,5).Value = 2.567607164382935
$ws.Cells.Item(2,6).Value = 0.3823793908064108
$ws.Cells.Item(2,7).Value = 0.1271628014236173
$ws.Cells.Item(2,8).Value = 0.7923793438529929

$ws.Cells.Item(3,3).Value = 0.8655490875244141
$ws.Cells.Item(3,4).Value = -18.99405670166016
$ws.Cells.Item(3,5).Value = 5.860441207885742
$ws.Cells.Item(3,6).Value = -0.654463132950649
$ws.Cells.Item(3,7).Value = 0.1985160146708309
$ws.Cells.Item(3,8).Value = 3.398013450712434

$ws.Cells.Item(4,3).Value = 12.034010887146
$ws.Cells.Item(4,4).Value = -13.59659099578857
$ws.Cells.Item(4,5).Value = 12.0263729095459
$ws.Cells.Item(4,6).Value = -0.3379178534003359
$ws.Cells.Item(4,7).Value = 0.9397114328064812
$ws.Cells.Item(4,8).Value = 2.251024748018254

$ws.Cells.Item(5,3).Value = 5.544707775115967
$ws.Cells.Item(5,4).Value = -11.79852485656738
$ws.Cells.Item(5,5).Value = -1.972949981689453
$ws.Cells.Item(5,6).Value = 0.9759882148647756
$ws.Cells.Item(5,7).Value = -0.4570092331052646
$ws.Cells.Item(5,8).Value = -1.31463850730377

$ws.Cells.Item(6,3).Value = 6.281956195831299
$ws.Cells.Item(6,4).Value = -19.22240829467773
$ws.Cells.Item(6,5).Value = -3.653380393981934
$ws.Cells.Item(6,6).Value = 0.7107792543491266
$ws.Cells.Item(6,7).Value = -2.305014370311596
$ws.Cells.Item(6,8).Value = -2.200793484742744

$ws.Cells.Item(7,3).Value = -1.975251197814941
$ws.Cells.Item(7,4).Value = -43.3292350769043
$ws.Cells.Item(7,5).Value = 0.271059513092041
$ws.Cells.Item(7,6).Value = 1.819485336698153
$ws.Cells.Item(7,7).Value = 0.8139906653559265
$ws.Cells.Item(7,8).Value = -4.443760702123214

$ws.Cells.Item(8,3).Value = -20.38084411621094
$ws.Cells.Item(8,4).Value = -15.68858814239502
$ws.Cells.Item(8,5).Value = 1.961378574371338
$ws.Cells.Item(8,6).Value = -0.05035778864519869
$ws.Cells.Item(8,7).Value = 1.772313159173635
$ws.Cells.Item(8,8).Value = -4.182391558642067

$ws.Cells.Item(9,3).Value = 10.76584243774414
$ws.Cells.Item(9,4).Value = -33.70746231079102
$ws.Cells.Item(9,5).Value = 31.38712692260743
$ws.Cells.Item(9,6).Value = -4.364394300271125
$ws.Cells.Item(9,7).Value = -2.634685883971363
$ws.Cells.Item(9,8).Value = 5.322240739592762

$ws.Cells.Item(10,3).Value = -27.7362232208252
$ws.Cells.Item(10,4).Value = -0.916855812072754
$ws.Cells.Item(10,5).Value = -10.06494903564453
$ws.Cells.Item(10,6).Value = -6.08671497237609
$ws.Cells.Item(10,7).Value = 2.640464704074256
$ws.Cells.Item(10,8).Value = 8.124276500721788

$ws.Cells.Item(11,3).Value = 8.33868408203125
$ws.Cells.Item(11,4).Value = -7.790350914001465
$ws.Cells.Item(11,5).Value = 11.65683746337891
$ws.Cells.Item(11,6).Value = 2.494929860399702
$ws.Cells.Item(11,7).Value = 7.262385212313916
$ws.Cells.Item(11,8).Value = 0.6411508452829229

$ws.Cells.Item(12,3).Value = -1.120648384094239
$ws.Cells.Item(12,4).Value = -23.34181785583496
$ws.Cells.Item(12,5).Value = -21.7900619506836
$ws.Cells.Item(12,6).Value = 5.578452347460839
$ws.Cells.Item(12,7).Value = -3.65870345949509
$ws.Cells.Item(12,8).Value = -2.740860200053031

$ws.Cells.Item(13,3).Value = 33.84098815917969
$ws.Cells.Item(13,4).Value = -33.77373504638672
$ws.Cells.Item(13,5).Value = -8.738304138183594
$ws.Cells.Item(13,6).Value = 5.107148140512828
$ws.Cells.Item(13,7).Value = -6.417216838576375
$ws.Cells.Item(13,8).Value = 0.5648952454172127

$ws.Cells.Item(14,3).Value = -9.054259300231934
$ws.Cells.Item(14,4).Value = -1.461036801338196
$ws.Cells.Item(14,5).Value = 6.970683097839356
$ws.Cells.Item(14,6).Value = 2.179000979318631
$ws.Cells.Item(14,7).Value = -0.9147634414166442
$ws.Cells.Item(14,8).Value = -6.292253988029406

$ws.Cells.Item(15,3).Value = 10.75043201446533
$ws.Cells.Item(15,4).Value = -19.06211471557617
$ws.Cells.Item(15,5).Value = 15.70715045928955
$ws.Cells.Item(15,6).Value = -4.380800012518604
$ws.Cells.Item(15,7).Value = -5.969521393950635
$ws.Cells.Item(15,8).Value = 0.8675613487578511

$ws.Cells.Item(16,3).Value = 1.280778884887695
$ws.Cells.Item(16,4).Value = -11.42607116699219
$ws.Cells.Item(16,5).Value = -8.684724807739258
$ws.Cells.Item(16,6).Value = -4.858158437369363
$ws.Cells.Item(16,7).Value = 3.663422577044154
$ws.Cells.Item(16,8).Value = 4.786025135304881

$ws.Cells.Item(17,3).Value = 36.76531219482422
$ws.Cells.Item(17,4).Value = -8.253963470458984
$ws.Cells.Item(17,5).Value = -22.78386306762696
$ws.Cells.Item(17,6).Value = -2.86569286016895
$ws.Cells.Item(17,7).Value = 9.76977250463675
$ws.Cells.Item(17,8).Value = -0.04881703791199454

$ws.Cells.Item(18,3).Value = -36.61545944213867
$ws.Cells.Item(18,4).Value = -18.09431838989257
$ws.Cells.Item(18,5).Value = 3.823569297790527
$ws.Cells.Item(18,6).Value = -3.519123984881036
$ws.Cells.Item(18,7).Value = 7.458226036650824
$ws.Cells.Item(18,8).Value = -5.801270187837286

$ws.Cells.Item(19,3).Value = 18.00795745849609
$ws.Cells.Item(19,4).Value = -45.09830856323242
$ws.Cells.Item(19,5).Value = -8.873518943786621
$ws.Cells.Item(19,6).Value = 2.29803731803496
$ws.Cells.Item(19,7).Value = -5.991696847046798
$ws.Cells.Item(19,8).Value = -5.515202137812292

$ws.Cells.Item(20,3).Value = -20.05809783935547
$ws.Cells.Item(20,4).Value = 0.3998336791992187
$ws.Cells.Item(20,5).Value = 3.845695495605469
$ws.Cells.Item(20,6).Value = 2.744500100300473
$ws.Cells.Item(20,7).Value = -0.9322929135791883
$ws.Cells.Item(20,8).Value = -6.360266283544435

$ws.Cells.Item(21,3).Value = 7.206372261047363
$ws.Cells.Item(21,4).Value = -20.33248519897461
$ws.Cells.Item(21,5).Value = 22.94344902038575
$ws.Cells.Item(21,6).Value = -5.222654512415387
$ws.Cells.Item(21,7).Value = -0.6178251758451783
$ws.Cells.Item(21,8).Value = 4.528362251701193

$ws.Cells.Item(22,1).Value = 2000
$ws.Cells.Item(22,2).Value = "walkingToRunning"
$ws.Cells.Item(22,3).Value = -10.1914873123169
$ws.Cells.Item(22,4).Value = -12.15236282348633
$ws.Cells.Item(22,5).Value = -10.96279335021973
$ws.Cells.Item(22,6).Value = -5.398936099406894
$ws.Cells.Item(22,7).Value = 0.6328901628238506
$ws.Cells.Item(22,8).Value = 6.700743228821124

$ws.Cells.Item(23,1).Value = 2100
$ws.Cells.Item(23,2).Value = "walkingToRunning"
$ws.Cells.Item(23,3).Value = 28.82599258422852
$ws.Cells.Item(23,4).Value = 6.167891502380371
$ws.Cells.Item(23,5).Value = -0.517308235168457
$ws.Cells.Item(23,6).Value = -2.420953338994997
$ws.Cells.Item(23,7).Value = 9.947797151136147
$ws.Cells.Item(23,8).Value = -0.9126796104522215

$ws.Cells.Item(24,1).Value = 2200
$ws.Cells.Item(24,2).Value = "walkingToRunning"
$ws.Cells.Item(24,3).Value = -22.66286087036133
$ws.Cells.Item(24,4).Value = -15.7267017364502
$ws.Cells.Item(24,5).Value = 0.2342269420623779
$ws.Cells.Item(24,6).Value = -0.05836680110209574
$ws.Cells.Item(24,7).Value = 8.935851176995778
$ws.Cells.Item(24,8).Value = -8.390108103527448

$ws.Cells.Item(25,1).Value = 2300
$ws.Cells.Item(25,2).Value = "walkingToRunning"
$ws.Cells.Item(25,3).Value = -11.30067539215088
$ws.Cells.Item(25,4).Value = -54.94432067871094
$ws.Cells.Item(25,5).Value = 17.55831718444824
$ws.Cells.Item(25,6).Value = 4.732191570142168
$ws.Cells.Item(25,7).Value = -8.935019577985051
$ws.Cells.Item(25,8).Value = -8.243645213661406

$ws.Cells.Item(26,1).Value = 2400
$ws.Cells.Item(26,2).Value = "walkingToRunning"
$ws.Cells.Item(26,3).Value = -17.29559326171875
$ws.Cells.Item(26,4).Value = 4.657787322998047
$ws.Cells.Item(26,5).Value = -2.186375617980957
$ws.Cells.Item(26,6).Value = 3.33527006372728
$ws.Cells.Item(26,7).Value = -10.98798226935694
$ws.Cells.Item(26,8).Value = -9.769438461483361

$ws.Cells.Item(27,1).Value = 2500
$ws.Cells.Item(27,2).Value = "walkingToRunning"
$ws.Cells.Item(27,3).Value = 20.14034080505371
$ws.Cells.Item(27,4).Value = -19.13811683654785
$ws.Cells.Item(27,5).Value = 23.86569976806641
$ws.Cells.Item(27,6).Value = -0.9925953642860486
$ws.Cells.Item(27,7).Value = -2.305789720325507
$ws.Cells.Item(27,8).Value = 3.166530378201873

$ws.Cells.Item(28,1).Value = 2600
$ws.Cells.Item(28,2).Value = "walkingToRunning"
$ws.Cells.Item(28,3).Value = 13.2857141494751
$ws.Cells.Item(28,4).Value = -19.34296989440918
$ws.Cells.Item(28,5).Value = -3.264841318130493
$ws.Cells.Item(28,6).Value = -4.912634885748004
$ws.Cells.Item(28,7).Value = 4.700535394758443
$ws.Cells.Item(28,8).Value = 3.522403849981169

$ws.Cells.Item(29,1).Value = 2700
$ws.Cells.Item(29,2).Value = "walkingToRunning"
$ws.Cells.Item(29,3).Value = -11.63338565826416
$ws.Cells.Item(29,4).Value = 5.394529819488525
$ws.Cells.Item(29,5).Value = -0.188831090927124
$ws.Cells.Item(29,6).Value = -3.346227371255714
$ws.Cells.Item(29,7).Value = 3.616383455186579
$ws.Cells.Item(29,8).Value = -1.352327122114249

$ws.Cells.Item(30,1).Value = 2800
$ws.Cells.Item(30,2).Value = "walkingToRunning"
$ws.Cells.Item(30,3).Value = -38.69764709472656
$ws.Cells.Item(30,4).Value = -13.45611763000488
$ws.Cells.Item(30,5).Value = 1.173340797424316
$ws.Cells.Item(30,6).Value = 3.400979183107231
$ws.Cells.Item(30,7).Value = 13.99426472249469
$ws.Cells.Item(30,8).Value = -1.426351852441958

$ws.Cells.Item(31,1).Value = 2900
$ws.Cells.Item(31,2).Value = "walkingToRunning"
$ws.Cells.Item(31,3).Value = 9.269144058227541
$ws.Cells.Item(31,4).Value = -44.79425811767578
$ws.Cells.Item(31,5).Value = 1.448012948036194
$ws.Cells.Item(31,6).Value = 0.4026054789882916
$ws.Cells.Item(31,7).Value = -6.858872019183412
$ws.Cells.Item(31,8).Value = -4.158537519539816
